$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "episode: 372000 win_rate: 0.305",
    "episode: 373000 win_rate: 0.332",
    "episode: 374000 win_rate: 0.322",
    "episode: 375000 win_rate: 0.3165",
    "episode: 376000 win_rate: 0.319",
    "episode: 377000 win_rate: 0.281",
    "episode: 378000 win_rate: 0.306",
    "episode: 379000 win_rate: 0.307",
    "episode: 380000 win_rate: 0.348",
    "episode: 381000 win_rate: 0.321",
    "episode: 382000 win_rate: 0.319",
    "episode: 383000 win_rate: 0.3",
    "episode: 384000 win_rate: 0.308",
    "episode: 385000 win_rate: 0.3375",
    "episode: 386000 win_rate: 0.324",
    "episode: 387000 win_rate: 0.3315",
    "episode: 388000 win_rate: 0.322",
    "episode: 389000 win_rate: 0.331",
    "episode: 390000 win_rate: 0.3105",
    "episode: 391000 win_rate: 0.322",
    "episode: 392000 win_rate: 0.323",
    "episode: 393000 win_rate: 0.298"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $srcCell = $ws.Cells.Item($row, 1)
    $dstCell = $ws.Cells.Item($row, 6)
    $dstCell.Value = $values[$i]
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null
}

$ws.Range("I23").Select()
